$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (RUG823.fasta) entirely.
$ws.Rows.Item(20).Delete()

# Remove the old "max" column (column C, constant value 1).
# This shifts old column D ("prediction") into C and old column E
# ("rejection-f") into D.
$ws.Columns.Item(3).Delete()

# Column B now holds the numeric prediction score per row instead of a
# constant 1.
$ws.Range("B2").Value = 285716.374927806
$ws.Range("B3").Value = -50347.13030853675
$ws.Range("B4").Value = 296501.3584019154
$ws.Range("B5").Value = 24940.4447828519
$ws.Range("B6").Value = 294066.3969261379
$ws.Range("B7").Value = 308236.8111055916
$ws.Range("B8").Value = -74547.32242456335
$ws.Range("B9").Value = 166118.2343003266
$ws.Range("B10").Value = 282873.4098115551
$ws.Range("B11").Value = 319705.5857858557
$ws.Range("B12").Value = 229580.9630112328
$ws.Range("B13").Value = 329823.4669816623
$ws.Range("B14").Value = 279222.8302084922
$ws.Range("B15").Value = 136400.7482931591
$ws.Range("B16").Value = 303175.0149353194
$ws.Range("B17").Value = 31656.89979797075
$ws.Range("B18").Value = 310390.5904991316
$ws.Range("B19").Value = 94639.06071850809

# Rows whose score was negative are marked as rejected in the
# rejection-f column (now column D).
$ws.Range("D3").Value = "o__Methanobacteriales(reject)"
$ws.Range("D8").Value = "o__Methanobacteriales(reject)"
